$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamps for column B (rows 2-30), generated by the RAD test-case
# runner on a later execution (Tue Jan 28 2025 run).
$timestamps = @(
    "Tue Jan 28 21:36:09 EST 2025",
    "Tue Jan 28 21:36:19 EST 2025",
    "Tue Jan 28 21:36:28 EST 2025",
    "Tue Jan 28 21:36:39 EST 2025",
    "Tue Jan 28 21:36:48 EST 2025",
    "Tue Jan 28 21:36:58 EST 2025",
    "Tue Jan 28 21:37:08 EST 2025",
    "Tue Jan 28 21:37:18 EST 2025",
    "Tue Jan 28 21:37:28 EST 2025",
    "Tue Jan 28 21:37:38 EST 2025",
    "Tue Jan 28 21:37:49 EST 2025",
    "Tue Jan 28 21:37:59 EST 2025",
    "Tue Jan 28 21:38:09 EST 2025",
    "Tue Jan 28 21:38:19 EST 2025",
    "Tue Jan 28 21:38:30 EST 2025",
    "Tue Jan 28 21:38:40 EST 2025",
    "Tue Jan 28 21:38:50 EST 2025",
    "Tue Jan 28 21:39:00 EST 2025",
    "Tue Jan 28 21:39:10 EST 2025",
    "Tue Jan 28 21:39:21 EST 2025",
    "Tue Jan 28 21:39:31 EST 2025",
    "Tue Jan 28 21:39:41 EST 2025",
    "Tue Jan 28 21:39:51 EST 2025",
    "Tue Jan 28 21:40:01 EST 2025",
    "Tue Jan 28 21:40:11 EST 2025",
    "Tue Jan 28 21:40:21 EST 2025",
    "Tue Jan 28 21:40:31 EST 2025",
    "Tue Jan 28 21:40:41 EST 2025",
    "Tue Jan 28 21:40:52 EST 2025"
)

# Rows whose Result flips from Pass to Fail in this run.
$failRows = @(18, 19, 29)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $timestamps[$i]
    if ($failRows -contains $row) {
        $ws.Cells.Item($row, 1).Value = "Fail"
    }
}
